$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (quarters 43465 and 43373),
# shifting the existing D:K data right to F:M.
$ws.Columns("D:E").Insert()

# Copy number formatting (date / comma-number styles) from column F
# into the two newly inserted blank columns.
$ws.Columns("F:F").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns with the new quarter data.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 9988600
$ws.Range("E8").Value = 10429500
$ws.Range("D9").Value = 9726400
$ws.Range("E9").Value = 10162800
$ws.Range("D10").Value = 262200
$ws.Range("E10").Value = 266700
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 17100
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 9925300
$ws.Range("E17").Value = 10351300
$ws.Range("D18").Value = 63300
$ws.Range("E18").Value = 78200
$ws.Range("D20").Value = -800
$ws.Range("E20").Value = 2700
$ws.Range("D21").Value = 85000
$ws.Range("E21").Value = 101500
$ws.Range("D22").Value = 19600
$ws.Range("E22").Value = 19100
$ws.Range("D23").Value = 42900
$ws.Range("E23").Value = 61800
$ws.Range("D24").Value = 13200
$ws.Range("E24").Value = 23000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 29700
$ws.Range("E26").Value = 38800
$ws.Range("D27").Value = 29600
$ws.Range("E27").Value = 38200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 800
$ws.Range("E32").Value = -2700
$ws.Range("D33").Value = 29600
$ws.Range("E33").Value = 38200
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 29600
$ws.Range("E35").Value = 38200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 211700
$ws.Range("E41").Value = 142100
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 2739600
$ws.Range("E43").Value = 3106900
$ws.Range("D44").Value = 523100
$ws.Range("E44").Value = 679300
$ws.Range("D45").Value = 500400
$ws.Range("E45").Value = 399600
$ws.Range("D46").Value = 3974800
$ws.Range("E46").Value = 4327900
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 350300
$ws.Range("E48").Value = 338800
$ws.Range("D49").Value = 852700
$ws.Range("E49").Value = 855600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 499100
$ws.Range("E52").Value = 485700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 5676900
$ws.Range("E54").Value = 6008000
$ws.Range("D57").Value = 2399600
$ws.Range("E57").Value = 2785900
$ws.Range("D58").Value = 41100
$ws.Range("E58").Value = 35700
$ws.Range("D59").Value = 495200
$ws.Range("E59").Value = 461800
$ws.Range("D60").Value = 2935900
$ws.Range("E60").Value = 3283400
$ws.Range("D61").Value = 659900
$ws.Range("E61").Value = 706600
$ws.Range("D62").Value = 249500
$ws.Range("E62").Value = 236800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 3861400
$ws.Range("E66").Value = 4242800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1606100
$ws.Range("E72").Value = 1580400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1815500
$ws.Range("E76").Value = 1765200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 29600
$ws.Range("E81").Value = 38200
$ws.Range("D83").Value = 22500
$ws.Range("E83").Value = 20600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 133500
$ws.Range("E89").Value = -87500
$ws.Range("D91").Value = -27600
$ws.Range("E91").Value = -15800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -13000
$ws.Range("E94").Value = 99500
$ws.Range("D96").Value = -4200
$ws.Range("E96").Value = -4000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -45100
$ws.Range("E100").Value = -57300
$ws.Range("D101").Value = -5700
$ws.Range("E101").Value = -200
$ws.Range("D102").Value = 69700
$ws.Range("E102").Value = -45500

# A handful of historical cells (now shifted into columns I/J) were
# corrected with restated figures as part of this update.
$ws.Range("I91").Value = -6300
$ws.Range("J91").Value = -21500
$ws.Range("H94").Value = 55300
$ws.Range("I94").Value = 67400
